$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp footer (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 04:05"

# Country stats refreshed / table re-sorted by "Casos totales" (rows shuffled
# because some totals changed). Each array is: row, country name, then the
# Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos /
# Muertes hoy / Muertes values that belong in columns B-H for that row.
$updates = @(
    ,@(67, "Bolivia", 4481, 218, 533, 3759, 0, 15, 189)
    ,@(74, "Honduras", 2955, 157, 349, 2459, 0, 1, 147)
    ,@(75, "Guinea", 2863, 0, 1525, 1320, 0, 0, 18)
    ,@(76, "Uzbekistan", 2855, 0, 2338, 504, 0, 0, 13)
    ,@(77, "Grecia", 2840, 0, 1374, 1301, 0, 0, 165)
    ,@(195, "Santa Lucia", 18, 0, 18, 0, 0, 0, 0)
    ,@(196, "Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)
    ,@(197, "Belice", 18, 0, 16, 0, 0, 0, 2)
    ,@(209, "Seychelles", 11, 0, 11, 0, 0, 0, 0)
    ,@(211, "Montserrat", 11, 0, 10, 0, 0, 0, 1)
    ,@(214, "Sahara Occidental", 6, 0, 6, 0, 0, 0, 0)
    ,@(216, "Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 1).Value = $u[1]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $u[$col]
    }
}
